# Update cryptos list (prices + 1h volume %) per the latest GitHub Actions scrape.
# Values that look like plain decimal numbers are written with a leading
# apostrophe so Excel stores them as text (matching the source data, which
# keeps trailing/insignificant-looking digits such as "67.60" or "0.112"
# intact instead of normalizing them as floating point numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.074.56'
$ws.Range("E2").Value = '  +0.58%  '
$ws.Range("D3").Value = '1.675.58'
$ws.Range("E3").Value = '  +0.35%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = "'215.22"
$ws.Range("E5").Value = '  +0.23%  '
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("E8").Value = '  +2.31%  '
$ws.Range("D9").Value = "'21.24"
$ws.Range("E9").Value = '  +4.72%  '
$ws.Range("D10").Value = "'0.0621"
$ws.Range("E10").Value = '  +0.13%  '
$ws.Range("E11").Value = '  -0.71%  '
$ws.Range("D12").Value = '1.912.49'
$ws.Range("E12").Value = '  +0.43%  '
$ws.Range("D13").Value = '1.673.64'
$ws.Range("E13").Value = '  +0.73%  '
$ws.Range("E14").Value = '  +0.93%  '
$ws.Range("D15").Value = "'0.536"
$ws.Range("E15").Value = '  +1.95%  '
$ws.Range("D16").Value = "'66.05"
$ws.Range("E16").Value = '  +0.90%  '
$ws.Range("D17").Value = '27.057.33'
$ws.Range("E17").Value = '  +0.55%  '
$ws.Range("D18").Value = "'237.05"
$ws.Range("E18").Value = '  +1.60%  '
$ws.Range("D19").Value = "'8.14"
$ws.Range("E19").Value = '  +1.59%  '
$ws.Range("D20").Value = '0.0₃0738'
$ws.Range("E20").Value = '  +0.67%  '
$ws.Range("E21").Value = '  +0.12%  '
$ws.Range("D22").Value = "'4.46"
$ws.Range("E22").Value = '  +1.13%  '
$ws.Range("D23").Value = "'9.29"
$ws.Range("E23").Value = '  +1.61%  '
$ws.Range("E24").Value = '  -1.46%  '
$ws.Range("D25").Value = "'147.30"
$ws.Range("E25").Value = '  +0.71%  '
$ws.Range("D26").Value = "'7.23"
$ws.Range("E26").Value = '  +1.68%  '
$ws.Range("D27").Value = "'16.31"
$ws.Range("E27").Value = '  +2.44%  '
$ws.Range("D28").Value = "'0.112"
$ws.Range("E28").Value = '  +0.52%  '
$ws.Range("E29").Value = '  +0.23%  '
$ws.Range("E30").Value = '  +0.12%  '
$ws.Range("E31").Value = '  +0.29%  '
$ws.Range("E32").Value = '  +0.71%  '
$ws.Range("D33").Value = '1.529.80'
$ws.Range("E33").Value = '  +5.01%  '
$ws.Range("D34").Value = "'3.17"
$ws.Range("E34").Value = '  +1.65%  '
$ws.Range("D35").Value = "'1.69"
$ws.Range("E35").Value = '  +3.27%  '
$ws.Range("E36").Value = '  -0.90%  '
$ws.Range("D37").Value = "'0.593"
$ws.Range("E37").Value = '  +1.58%  '
$ws.Range("D38").Value = "'0.915"
$ws.Range("E38").Value = '  +1.76%  '
$ws.Range("E39").Value = '  +2.42%  '
$ws.Range("D40").Value = "'1.07"
$ws.Range("E40").Value = '  +2.58%  '
$ws.Range("E41").Value = '  +0.13%  '
$ws.Range("D42").Value = "'67.60"
$ws.Range("E42").Value = '  +2.10%  '
$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").Value = "'2.26"
$ws.Range("E43").Value = '  -1.32%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = "'5.51"
$ws.Range("E44").Value = '  -3.83%  '
$ws.Range("D45").Value = '1.820.06'
$ws.Range("E45").Value = '  +0.57%  '
$ws.Range("E46").Value = '  +0.42%  '
$ws.Range("D47").Value = "'90.67"
$ws.Range("E47").Value = '  -0.06%  '
$ws.Range("E48").Value = '  +1.05%  '
$ws.Range("E49").Value = '  +2.36%  '
$ws.Range("D50").Value = "'7.97"
$ws.Range("E50").Value = '  +5.05%  '
$ws.Range("E51").Value = '  +0.45%  '
